{"js": "// Fixed #418 Empty AQL expressions generate empty lines.\n// Remove the empty paragraph (no visible text, just leftover run formatting\n// such as the orange font color) that an empty AQL expression left behind.\n\nconst body = context.document.body;\nconst paragraphs = body.paragraphs;\nparagraphs.load(\"items\");\nawait context.sync();\n\n// Load text and font color for every paragraph so we can find the\n// specific empty paragraph to drop (there can be more than one empty\n// paragraph in the document, e.g. a trailing one at the very end that\n// must be preserved).\nparagraphs.items.forEach((p) => {\n  p.load(\"text\");\n  p.font.load(\"color\");\n});\nawait context.sync();\n\nlet target = null;\nfor (const p of paragraphs.items) {\n  const isEmpty = p.text === \"\" || p.text === \"\\r\";\n  if (isEmpty && p.font.color && p.font.color.toUpperCase() === \"#E36C0A\") {\n    target = p;\n    break;\n  }\n}\n\n// Fallback: if no paragraph matched the color (should not normally\n// happen), fall back to the first empty paragraph that is not the very\n// last paragraph of the body (that last one is intentionally kept).\nif (!target) {\n  const items = paragraphs.items;\n  for (let i = 0; i < items.length - 1; i++) {\n    const p = items[i];\n    if (p.text === \"\" || p.text === \"\\r\") {\n      target = p;\n      break;\n    }\n  }\n}\n\nif (target) {\n  target.delete();\n  await context.sync();\n}\n", "ps1": "# Fixed #418 Empty AQL expressions generate empty lines.\n# Remove the empty paragraph (no visible text, just leftover run formatting\n# such as the orange font color) that an empty AQL expression left behind.\n# The very last paragraph of the document is intentionally left untouched.\n\n$d = $word.ActiveDocument\n$count = $d.Paragraphs.Count\n\n$target = $null\nfor ($i = 1; $i -le $count; $i++) {\n    if ($i -eq $count) {\n        # never remove the final paragraph of the document\n        continue\n    }\n    $p = $d.Paragraphs.Item($i)\n    $r = $p.Range\n    # An \"empty\" paragraph only contains its terminating paragraph mark,\n    # so Range.Text has length 1 (just \"\\r\").\n    if ($r.Text.Length -eq 1) {\n        $target = $p\n        break\n    }\n}\n\nif ($target -ne $null) {\n    $target.Range.Delete()\n}\n"}
